# elapsed time y cpu
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Data rows 2-14
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 7).Value = 1.669922641383406
    $ws.Cells.Item($row, 8).Value = 0.97
}
